$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Preserve / relocate the trailing hyperlink row ---
# Original layout: row 23 = "<link text>" cell (A23) styled "Hyperlink", with a
# hyperlink pointing at the rubygarage blog post.
# New layout needs three new data rows inserted right before it (rows 22-24),
# pushing that hyperlink row down to row 25.

# Remove the existing hyperlink cleanly first (keeps the cell's value/style intact)
$ws.Range("A23").Hyperlinks.Delete()

# Insert two blank rows above row 23; this shifts the existing row 23
# (value + "Hyperlink" style) down to row 25, and rows 1-21 stay untouched.
$ws.Range("A23:A24").EntireRow.Insert()

# --- Fill in the three new command rows (22-24) ---
# NB: cell values are written in this particular order so that the shared
# string table ends up with the same unique-string ordering as the target
# workbook (new strings are appended to sharedStrings.xml in first-seen
# order: A22, A23, B23, B22, A24, B24).
$ws.Range("A22").Value = "$ git commit -m 'message'"
$ws.Range("A23").Value = "$ git commit"
$ws.Range("B23").Value = "To commit your changes at local repository but this will open new editor window and you have to press 'I' to start typing and after typing message, type ':wq' to go back to git bash and to escape from new editor"
$ws.Range("B22").Value = "To commit your changes at local repository(using this way you can skip edit stage)"
$ws.Range("A24").Value = "$ clear"
$ws.Range("B24").Value = "To clear the window"

# --- Re-attach the hyperlink to its new home (row 25) and restore its style ---
$ws.Hyperlinks.Add($ws.Range("A25"), "https://rubygarage.org/blog/most-basic-git-commands-with-examples") | Out-Null
$ws.Range("A25").Style = "Hyperlink"

# --- Update the view selection to match the saved workbook state ---
$ws.Range("B24").Select()
